$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 is the "祭坛" (Altar) card. Its effect text (column C) changes the
# damage cost of option 3 from 3 to 2, so the card stays longer in the
# player's Loot Zone: "受到3伤害" (take 3 damage) -> "受到2伤害" (take 2 damage).
$newText = "多选：①将1张手牌送墓，翻开遭遇牌堆顶1张牌，如果是战利品牌则可以获得。②弃置1张战利品牌，获得1道具点。③受到2伤害，获得1SP。"

$ws.Range("C11").Value = $newText
